$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new column K that extends the existing table (years 2014-2020 in D:J)
# with the 2021 data column, copying the formatting from column J row-by-row.

# Row 2: trailing border-only cell under the header, no value, same style as J2.
$ws.Range("J2").Copy()
$ws.Range("K2").PasteSpecial(-4122)

# Row 3: year header value 2021, same style as J3.
$ws.Range("J3").Copy()
$ws.Range("K3").PasteSpecial(-4122)
$ws.Range("K3").Value = 2021

# Row 4
$ws.Range("J4").Copy()
$ws.Range("K4").PasteSpecial(-4122)
$ws.Range("K4").Value = 295

# Row 5
$ws.Range("J5").Copy()
$ws.Range("K5").PasteSpecial(-4122)
$ws.Range("K5").Value = 163

# Row 6
$ws.Range("J6").Copy()
$ws.Range("K6").PasteSpecial(-4122)
$ws.Range("K6").Value = 268

# Row 7
$ws.Range("J7").Copy()
$ws.Range("K7").PasteSpecial(-4122)
$ws.Range("K7").Value = 155

# Row 8
$ws.Range("J8").Copy()
$ws.Range("K8").PasteSpecial(-4122)
$ws.Range("K8").Value = 27

# Row 9
$ws.Range("J9").Copy()
$ws.Range("K9").PasteSpecial(-4122)
$ws.Range("K9").Value = 8

$excel.CutCopyMode = 0

# Update the saved selection to match the new active cell in the workbook.
[void]$ws.Range("L5").Select()
